# Edit script for Pedidos.xlsx
# Replaces the order list data (rows 2-63 of column A:C) on the active sheet
# with an updated, much shorter list of shipments, clearing the rows that no
# longer have data, and moves the active selection to A19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (Remessa, Material, Quantidade) replacing the old rows 2-11.
$data = @(
    @("80264004", "33569-ATE-I", 1),
    @("80265160", "30018-KMT-I", 24000),
    @("80265942", "20637-TDK-I", 2000),
    @("80265942", "20850-FUZ-I", 600),
    @("80265942", "20869-FUZ-I", 36000),
    @("80265942", "20853-FUZ-I", 1900),
    @("80265942", "60234-WUE-I", 1500),
    @("80265944", "20389-DCC-I", 5000),
    @("80266324", "40193-TDK-N", 15000),
    @("80266631", "12350-DLO-I", 1)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Remaining previously-filled rows (12-63) no longer have data; clear their
# contents so the cells become blank while keeping their existing formatting.
$lastOldRow = 63
$firstClearRow = $startRow + $data.Count
if ($firstClearRow -le $lastOldRow) {
    $clearRange = $ws.Range("A" + $firstClearRow + ":C" + $lastOldRow)
    $clearRange.ClearContents()
}

# Move the active selection to A19, matching the saved sheet view state.
$ws.Range("A19").Select()
